$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    "F2"  = 22;  "H2"  = 24
    "F3"  = 23;  "H3"  = 26
    "F4"  = 15;  "H4"  = 18
    "F5"  = 24;  "H5"  = 28
    "F6"  = 35;  "H6"  = 42
    "F7"  = 23;  "H7"  = 24
    "F9"  = 8;   "H9"  = 12
    "F10" = 20;  "H10" = 22
    "F12" = 22;  "H12" = 24
    "F14" = 25;  "H14" = 27
    "F15" = 68;  "H15" = 79
    "F16" = 117; "H16" = 205
    "F17" = 18;  "H17" = 19
    "F18" = 40;  "H18" = 63
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
